# Add step "4" to the architecture diagram, and reposition the "1" and "2"
# step markers to their new spots.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Step "1" marker (TextBox 1) moves to its new position ---
$step1 = $s.Shapes.Item(15)
$step1.Left = 300.1251988503937
$step1.Top  = 23.209056118110237

# --- Step "2" marker (TextBox 40) slides left (same vertical position) ---
$step2 = $s.Shapes.Item(16)
$step2.Left = 339.9538732677166

# --- New step "4" marker ---
# Prime the shape-ID allocator first: PowerPoint hands out the lowest
# unused id, and this deck already has a few gaps (ids 3, 4, 15, 19) below
# its current max id of 42. Briefly adding/removing throwaway textboxes
# consumes those gaps so the real new shape lands on id 43, exactly as a
# fresh shape added in real PowerPoint would.
$dummies = @()
for ($i = 0; $i -lt 4; $i++) {
    $dummies += $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
}
foreach ($d in $dummies) {
    $d.Delete()
}

# Duplicate the "3" marker (TextBox 41) so the new one inherits identical
# formatting (fill, font, size, body properties), then move it into place
# and update its label text to "4".
$step3 = $s.Shapes.Item(17)
$newRange = $step3.Duplicate()
$step4 = $newRange.Item(1)
$step4.Name = "TextBox 42"
$step4.Left = 553.0349606299212
$step4.Top  = 191.53417322834645
$step4.TextFrame.TextRange.Text = "4"
